$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F header: United States
$ws.Range("F1").Value = "United States"

# Row 2 (Global tax on millionaires)
$ws.Range("B2").Value = 0.793548387096774
$ws.Range("C2").Value = 0.9609375
$ws.Range("D2").Value = 1.20754716981132
$ws.Range("E2").Value = 0.8801261829653
$ws.Range("F2").Font.Bold = $false

# Row 3 (National tax on millionaires)
$ws.Range("B3").Value = 0.751612903225806
$ws.Range("C3").Value = 1.0703125
$ws.Range("D3").Value = 1.22264150943396
$ws.Range("E3").Value = 1.11356466876972
$ws.Range("F3").Font.Bold = $false

# Row 4 (Payments from high-income countries...)
$ws.Range("B4").Value = 0.306451612903226
$ws.Range("C4").Value = 0.359313077939234
$ws.Range("D4").Value = 0.756906077348066
$ws.Range("E4").Value = 0.277950310559006
$ws.Range("F4").Value = 0.0408626560726447

# Row 5 (High-income countries funding renewable energy...)
$ws.Range("B5").Value = 0.52258064516129
$ws.Range("C5").Value = 0.696169088507266
$ws.Range("D5").Value = 0.869244935543278
$ws.Range("E5").Value = 0.631987577639752
$ws.Range("F5").Value = 0.348467650397276

# Row 6 ($100 billion per year...)
$ws.Range("B6").Value = 0.403225806451613
$ws.Range("C6").Value = 0.578599735799207
$ws.Range("D6").Value = 0.69060773480663
$ws.Range("E6").Value = 0.389751552795031
$ws.Range("F6").Value = 0.149829738933031

# Row 7 (Cancellation of low-income countries' public debt)
$ws.Range("B7").Value = 0.0596774193548387
$ws.Range("C7").Value = -0.151915455746367
$ws.Range("D7").Value = 0.204419889502762
$ws.Range("E7").Value = 0.125776397515528
$ws.Range("F7").Value = -0.121831252364737

# Row 8 (Democratise international institutions...)
$ws.Range("B8").Value = 0.246774193548387
$ws.Range("C8").Value = 0.264200792602378
$ws.Range("D8").Value = 0.558011049723757
$ws.Range("E8").Value = 0.301242236024845
$ws.Range("F8").Value = 0.0620506999621642

# Row 9 (Removing tariffs on imports...)
$ws.Range("B9").Value = 0.154838709677419
$ws.Range("C9").Value = 0.408190224570674
$ws.Range("D9").Value = 0.548802946593002
$ws.Range("E9").Value = 0.607142857142857
$ws.Range("F9").Value = 0.178584941354521

# Row 10 (A minimum wage in all countries...)
$ws.Range("B10").Value = 0.479032258064516
$ws.Range("C10").Value = 0.513870541611625
$ws.Range("D10").Value = 0.710865561694291
$ws.Range("E10").Value = 0.515527950310559
$ws.Range("F10").Value = 0.195611048051457

# Row 11 (Fight tax evasion...)
$ws.Range("B11").Value = 1.00806451612903
$ws.Range("C11").Value = 0.895640686922061
$ws.Range("D11").Value = 1.11786372007366
$ws.Range("E11").Value = 0.992236024844721
$ws.Range("F11").Value = 0.184638668180098

# Row 12 (A maximum wealth limit...)
$ws.Range("B12").Value = 0.141935483870968
$ws.Range("C12").Value = 0.235138705416116
$ws.Range("D12").Value = 0.270718232044199
$ws.Range("E12").Value = 0.411490683229814
$ws.Range("F12").Value = -0.11312902005297
